$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in column D (T3) grades for rows 2-7
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1.5
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 2.5

# Update the active selection to D8 (matches sheetView selection in the diff)
$ws.Range("D8").Select()

$wb.Save()
